# Generate Report for Handoff
# Adds a new file (1bb5dfa7-4e2a-4023-a067-b3de17122f63.md) as row 3 to each of
# the three report sheets (Overview, zh-cn, de-de), mirroring the existing
# row 2 (0c838229-...) layout/format.

$wb = $excel.ActiveWorkbook

$newId   = "1bb5dfa7-4e2a-4023-a067-b3de17122f63"
$newMd   = "$newId.md"
$newPath = "e2e\$newId.md"
$newUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c54ffc605c34707458bc3ff8bdf88f220908f72f/e2e/$newId.md"

$zhXlf   = "$newId.8dea081ecd57f9e347f3c4534f2ec65e476d5f2c.zh-cn.xlf"
$deXlf   = "$newId.8dea081ecd57f9e347f3c4534f2ec65e476d5f2c.de-de.xlf"

$dtOverview = "2016-08-31 18:46:49"
$dtZh       = "2016-08-31 18:46:44"
$dtDe       = "2016-08-31 18:46:49"
$dtFmt      = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Range("B3").Value = $newPath
$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", $newPath) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $dtOverview
$wsOverview.Range("G3").NumberFormat = $dtFmt

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newMd
$wsZh.Range("A3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newUrl, "", "", $newMd) | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $dtZh
$wsZh.Range("H3").NumberFormat = $dtFmt
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = $dtFmt
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newMd
$wsDe.Range("A3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newUrl, "", "", $newMd) | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $dtDe
$wsDe.Range("H3").NumberFormat = $dtFmt
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = $dtFmt
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

# ---------------------------------------------------------------------------
# Grow the tables / autofilters so the new row is included.
# ---------------------------------------------------------------------------
foreach ($lo in $wsOverview.ListObjects) { $lo.Resize($wsOverview.Range("A1:G3")) }
foreach ($lo in $wsZh.ListObjects)       { $lo.Resize($wsZh.Range("A1:P3")) }
foreach ($lo in $wsDe.ListObjects)       { $lo.Resize($wsDe.Range("A1:P3")) }
